$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: add "Active" column (H1), copying the style from G1 ---
$ws.Range("H1").Value = "Active"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# --- Row 12 / Row 13: replace the old ". . . ." placeholder text ---
# (F13 is set first so the new shared strings land in the same table order as the target file)
$ws.Range("F13").Value = "나래이션 양식"
$ws.Range("F12").Value = "독백, 방백 양식"

# --- Row 13: D13 / E13 are no longer used (Character / State cleared entirely) ---
$ws.Range("D13:E13").Clear()

# --- New "Active" boolean flags (FALSE) on rows 8, 10, 11, 12, copying style from column G ---
$ws.Range("H8").Value = $false
$ws.Range("G8").Copy()
$ws.Range("H8").PasteSpecial(-4122)

$ws.Range("H10").Value = $false
$ws.Range("G10").Copy()
$ws.Range("H10").PasteSpecial(-4122)

$ws.Range("H11").Value = $false
$ws.Range("G11").Copy()
$ws.Range("H11").PasteSpecial(-4122)

$ws.Range("H12").Value = $false
$ws.Range("G12").Copy()
$ws.Range("H12").PasteSpecial(-4122)

# --- Selection / view state ---
$ws.Range("H2").Select()
